$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.945.78"
$ws.Range("E2").Value = "  -1.70%  "

$ws.Range("D3").Value = "3.130.51"
$ws.Range("E3").Value = "  -7.88%  "

$ws.Range("D5").Value = "569.05"
$ws.Range("E5").Value = "  -2.21%  "

$ws.Range("D6").Value = "168.50"
$ws.Range("E6").Value = "  -6.13%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  -3.15%  "

$ws.Range("D9").Value = "3.132.22"
$ws.Range("E9").Value = "  -7.78%  "

$ws.Range("E10").Value = "  -5.81%  "

$ws.Range("D11").Value = "6.52"
$ws.Range("E11").Value = "  -5.97%  "

$ws.Range("E12").Value = "  -5.92%  "

$ws.Range("D13").Value = "3.675.89"
$ws.Range("E13").Value = "  -7.85%  "

$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").Value = "26.61"
$ws.Range("E15").Value = "  -8.27%  "

$ws.Range("D16").Value = "64.833.63"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("E17").Value = "  -6.21%  "

$ws.Range("D18").Value = "3.129.79"
$ws.Range("E18").Value = "  -8.11%  "

$ws.Range("D19").Value = "5.67"
$ws.Range("E19").Value = "  -3.46%  "

$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  -7.06%  "

$ws.Range("D21").Value = "354.27"
$ws.Range("E21").Value = "  -3.30%  "

$ws.Range("D22").Value = "7.20"
$ws.Range("E22").Value = "  -4.50%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "68.87"

$ws.Range("D25").Value = "0.491"
$ws.Range("E25").Value = "  -7.11%  "

$ws.Range("D26").Value = "3.267.09"
$ws.Range("E26").Value = "  -7.83%  "

$ws.Range("E27").Value = "  -8.21%  "

$ws.Range("D28").Value = "9.60"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("E29").Value = "  -2.42%  "

$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("E32").Value = "  -4.13%  "

$ws.Range("D33").Value = "21.71"
$ws.Range("E33").Value = "  -6.32%  "

$ws.Range("D34").Value = "5.22"
$ws.Range("E34").Value = "  -8.86%  "

$ws.Range("E35").Value = "  -6.34%  "

$ws.Range("E36").Value = "  -5.73%  "

$ws.Range("D37").Value = "158.54"
$ws.Range("E37").Value = "  -1.82%  "

$ws.Range("E38").Value = "  -6.65%  "

$ws.Range("D39").Value = "0.825"
$ws.Range("E39").Value = "  -3.73%  "

$ws.Range("D40").Value = "26.11"
$ws.Range("E40").Value = "  -3.75%  "

$ws.Range("D41").Value = "1.74"
$ws.Range("E41").Value = "  -1.62%  "

$ws.Range("D42").Value = "2.646.10"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "6.07"
$ws.Range("E43").Value = "  -2.31%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.40"
$ws.Range("E44").Value = "  -7.65%  "

$ws.Range("E45").Value = "  -4.71%  "

$ws.Range("D46").Value = "39.34"
$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("D47").Value = "0.0649"
$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").Value = "23.78"
$ws.Range("E48").Value = "  -3.03%  "

$ws.Range("D49").Value = "316.60"
$ws.Range("E49").Value = "  -5.78%  "

$ws.Range("E50").Value = "  -5.38%  "

$ws.Range("E51").Value = "  -1.99%  "
